$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: advance the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value2 = 45309

# Step 2: update the prices in D23:D26 to 844
$ws.Range("D23").Value2 = 844
$ws.Range("D24").Value2 = 844
$ws.Range("D25").Value2 = 844
$ws.Range("D26").Value2 = 844
